# Applies the "tdwi.docx" edit: consolidates a handful of runs that were
# only split apart by spell/grammar-check proofing marks, and appends a
# large block of new content (web interface / usage / contact / thanks)
# right after the existing "future:" / cold-start paragraphs.

$d = $word.ActiveDocument

# --- 1. "...hosted by deloitte and tdwi. For future..." -------------------
# Runs "deloitte" | " and " | "tdwi" | ". For future" collapse into
# "deloitte" | " and tdwi. For future"
$d.Content.Find.Execute(" and tdwi. For future", $true, $false, $false, $false, $false, $true, 1, $false, " and tdwi. For future", 2) | Out-Null

# --- 2. "data used:" (was "data" + " used:" around a gramStart/gramEnd) ---
$d.Content.Find.Execute("data used:", $true, $false, $false, $false, $false, $true, 1, $false, "data used:", 2) | Out-Null

# --- 3. "I did not uses tripadvisor data because..." -----------------------
$d.Content.Find.Execute("I did not uses tripadvisor data because as stated on their homepage the API is not intended for data scrapping.", $true, $false, $false, $false, $false, $true, 1, $false, "I did not uses tripadvisor data because as stated on their homepage the API is not intended for data scrapping.", 2) | Out-Null

# --- 4. "(size under 100mb) ... noSQL database like graph based neo4j. ..." -
$d.Content.Find.Execute("(size under 100mb) I did not see the benefit of using a noSQL database like graph based neo4j.", $true, $false, $false, $false, $false, $true, 1, $false, "(size under 100mb) I did not see the benefit of using a noSQL database like graph based neo4j.", 2) | Out-Null

# --- 5. "nlp:" (was "nlp" + ":") -------------------------------------------
$d.Content.Find.Execute("nlp:", $true, $false, $false, $false, $false, $true, 1, $false, "nlp:", 2) | Out-Null

# --- 6. the long NLP paragraph (tf_idf / numpy spell-check splits) --------
$d.Content.Find.Execute("term frequency - inverse document frequency (tf_idf) and latent semantic analysis", $true, $false, $false, $false, $false, $true, 1, $false, "term frequency - inverse document frequency (tf_idf) and latent semantic analysis", 2) | Out-Null
$d.Content.Find.Execute("The result of the tf_idf is used in the LSA model.", $true, $false, $false, $false, $false, $true, 1, $false, "The result of the tf_idf is used in the LSA model.", 2) | Out-Null
$d.Content.Find.Execute("The importance of words (output of tf_idf) is used in the concept building", $true, $false, $false, $false, $false, $true, 1, $false, "The importance of words (output of tf_idf) is used in the concept building", 2) | Out-Null
$d.Content.Find.Execute("rows. I used numpy to calculate", $true, $false, $false, $false, $false, $true, 1, $false, "rows. I used numpy to calculate", 2) | Out-Null

# --- 7. "...what kind of service the datasubject values in each type..." --
$d.Content.Find.Execute("what kind of service the datasubject values in each type. I just gave it my best guess.", $true, $false, $false, $false, $false, $true, 1, $false, "what kind of service the datasubject values in each type. I just gave it my best guess.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Append the new closing section right after the "cold start" paragraph,
# which is currently the very last paragraph in the document.
# ---------------------------------------------------------------------------
$last = $d.Paragraphs($d.Paragraphs.Count)
$r = $last.Range
$r.Collapse(0)

$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("web interface:")
$r.Collapse(0)

$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("I’m not familiar with flask so I did not use the provided app for the web interface. But it should be easily possible to GET/POST the provided data.")
$r.Collapse(0)

$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("usage:")
$r.Collapse(0)

$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("you just need to create a virtual environment and install the req.txt, go into the main folder (travel_recommender) and run the app.py.")
$r.Collapse(0)

$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("You can reach me at:")
$r.Collapse(0)

$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("Benjamin Pohl")
$r.Collapse(0)
$r.InsertAfter("`vBenjamin.pohl95@gmail.com ")
$r.Collapse(0)

$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("Thank you for the challenge ")
$r.Collapse(0)

Write-Host "done"
